# Feb 15th #4 commit
# Update the "Final Result" status cells on each sheet with the latest
# QA test-run timestamps.

$wb = $excel.ActiveWorkbook

$wsSchoolSearch   = $wb.Worksheets.Item("School Search")
$wsProductSearch  = $wb.Worksheets.Item("Product Search")
$wsShoppingCart   = $wb.Worksheets.Item("Shopping Cart")
$wsCheckout       = $wb.Worksheets.Item("Checkout")
$wsPayment        = $wb.Worksheets.Item("Payment")

# School Search
$wsSchoolSearch.Range("C2").Value = "Success - 2021/02/15 20:47:57"
$wsSchoolSearch.Range("C3").Value = "Success - 2021/02/15 20:48:00"

# Product Search
$wsProductSearch.Range("K1").Value = "Success - 2021/02/15 20:49:27"

# Shopping Cart
$wsShoppingCart.Range("G2").Value = "Success - 2021/02/15 20:49:30"
$wsShoppingCart.Range("G3").Value = "Success - 2021/02/15 20:49:30"
$wsShoppingCart.Range("G4").Value = "Success - 2021/02/15 20:49:30"

# Checkout
$wsCheckout.Range("P2").Value = "Success - 2021/02/15 20:49:39"
$wsCheckout.Range("P3").Value = "Success - 2021/02/15 20:49:50"
$wsCheckout.Range("P4").Value = "Success - 2021/02/15 20:49:58"

# Payment
$wsPayment.Range("F2").Value = "Success - 2021/02/15 20:50:14"
